$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated pasture (Загон) statistics for rows 2-8, columns B:F
# (Сумма, Cреднаяя, Медианная, Макс, Мин)

$ws.Range("B2").Value = 0.0241
$ws.Range("C2").Value = 0.0241
$ws.Range("D2").Value = 0.0328
$ws.Range("E2").Value = 0.1438
$ws.Range("F2").Value = -0.1177

$ws.Range("B3").Value = 0.0106
$ws.Range("C3").Value = 0.0106
$ws.Range("D3").Value = 0.0157
$ws.Range("E3").Value = 0.1964
$ws.Range("F3").Value = -0.3781

$ws.Range("B4").Value = 0.0073
$ws.Range("C4").Value = 0.0073
$ws.Range("D4").Value = 0.0054
$ws.Range("E4").Value = 0.1344
$ws.Range("F4").Value = -0.1247

$ws.Range("B5").Value = 0.0083
$ws.Range("C5").Value = 0.0083
$ws.Range("D5").Value = 0.0145
$ws.Range("E5").Value = 0.1099
$ws.Range("F5").Value = -0.1902

$ws.Range("B6").Value = -0.0587
$ws.Range("C6").Value = -0.0587
$ws.Range("D6").Value = -0.0637
$ws.Range("E6").Value = 0.0508
$ws.Range("F6").Value = -0.1933

$ws.Range("B7").Value = -0.0576
$ws.Range("C7").Value = -0.0576
$ws.Range("D7").Value = -0.06
$ws.Range("E7").Value = 0.0753
$ws.Range("F7").Value = -0.1859

$ws.Range("B8").Value = 0.0433
$ws.Range("C8").Value = 0.0433
$ws.Range("D8").Value = 0.0461
$ws.Range("E8").Value = 0.1274
$ws.Range("F8").Value = -0.1028
